# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.692.96'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.937.20'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.548'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.72%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -5.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.97%  '
$ws.Range('E11').Value = '  +2.21%  '
$ws.Range('E12').Value = '  -3.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.47%  '
$ws.Range('D14').Value = '3.402.00'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('E15').Value = '  -5.74%  '
$ws.Range('D16').Value = '2.935.71'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.977'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '51.647.72'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.43%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.11%  '
$ws.Range('E26').Value = '  -6.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.108'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('E32').Value = '  -5.51%  '
$ws.Range('E33').Value = '  -5.46%  '
$ws.Range('E34').Value = '  -7.25%  '
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('E36').Value = '  -4.43%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.37%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.22'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.38%  '
$ws.Range('E40').Value = '  -6.43%  '
$ws.Range('E41').Value = '  -5.58%  '
$ws.Range('E42').Value = '  -4.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '120.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').Value = '2.095.77'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('E47').Value = '  -7.67%  '
$ws.Range('E48').Value = '  -6.69%  '
$ws.Range('D49').Value = '3.230.21'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  -5.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0318'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.05%  '
